$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for the third exam (Prova 03), copying the formatting from C6
$ws.Range("D6").Value = "Prova 03"
$ws.Range("C6").Copy()
$ws.Range("D6").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the "Prova 03" grades for every student row
$ws.Range("D7").Formula = "=15+15+15+20+15+15"
$ws.Range("D8").Formula = "=10+15+0+20+9+10"
$ws.Range("D9").Value = 0
$ws.Range("D10").Formula = "=3+10+15+20+12+0"
$ws.Range("D11").Value = 0
$ws.Range("D12").Formula = "=15+15+15+0+15+15"
$ws.Range("D13").Formula = "=10+15+15+20+15+20"
$ws.Range("D14").Formula = "=15+15+15+20+15+15"
$ws.Range("D15").Formula = "=15+10+15+20+12+0"
$ws.Range("D16").Formula = "=15+15+15+20+15+12"
$ws.Range("D17").Formula = "=15+5+15+20+12+15"
$ws.Range("D18").Value = 100
$ws.Range("D19").Value = 100
$ws.Range("D20").Value = 0
$ws.Range("D21").Formula = "=15+15+15+20+15+18"
$ws.Range("D22").Formula = "=15+15+15+20+15+15"
$ws.Range("D23").Formula = "=15+15+15+20+15+15"
$ws.Range("D24").Value = 0
$ws.Range("D25").Formula = "=8+12+15+20+12+15"
$ws.Range("D26").Formula = "=15+15+15+20+15+15"
$ws.Range("D27").Formula = "=15+15+15+20+15+10"
$ws.Range("D28").Formula = "=15+15+15+20+15+10"
$ws.Range("D29").Formula = "=15+12+15+20+12+12"
$ws.Range("D30").Formula = "=15+15+15+20+15+15"

# Match the final cell selection recorded in the workbook
$ws.Range("F8").Select()
